$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to retain text formatting so purely-numeric-looking
# strings (e.g. "0.998") are not auto-converted to numbers by Excel, matching
# the source workbook's inline-string cell type.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '34.119.85'
$ws.Range("E2").Value = '  +11.68%  '
$ws.Range("D3").Value = '1.820.23'
$ws.Range("E3").Value = '  +8.55%  '
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").Value = '228.64'
$ws.Range("E5").Value = '  +3.93%  '
$ws.Range("D6").Value = '0.547'
$ws.Range("E6").Value = '  +3.18%  '
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").Value = '31.42'
$ws.Range("E8").Value = '  +4.82%  '
$ws.Range("D9").Value = '47.35'
$ws.Range("E9").Value = '  +6.92%  '
$ws.Range("E10").Value = '  +6.45%  '
$ws.Range("D11").Value = '0.0668'
$ws.Range("E11").Value = '  +4.85%  '
$ws.Range("E12").Value = '  +2.47%  '
$ws.Range("D13").Value = '2.080.25'
$ws.Range("E13").Value = '  +8.42%  '
$ws.Range("D14").Value = '1.807.68'
$ws.Range("E14").Value = '  +7.54%  '
$ws.Range("D15").Value = '0.645'
$ws.Range("E15").Value = '  +4.54%  '
$ws.Range("D16").Value = '10.40'
$ws.Range("E16").Value = '  +1.28%  '
$ws.Range("D17").Value = '34.080.75'
$ws.Range("E17").Value = '  +11.46%  '
$ws.Range("D18").Value = '4.28'
$ws.Range("E18").Value = '  +7.68%  '
$ws.Range("D19").Value = '69.51'
$ws.Range("E19").Value = '  +4.76%  '
$ws.Range("D20").Value = '258.38'
$ws.Range("E20").Value = '  +5.47%  '
$ws.Range("D21").Value = '0.0₃0751'
$ws.Range("E21").Value = '  +3.76%  '
$ws.Range("D22").Value = '0.999'
$ws.Range("D23").Value = '10.56'
$ws.Range("E23").Value = '  +5.33%  '
$ws.Range("D24").Value = '4.34'
$ws.Range("E24").Value = '  +1.71%  '
$ws.Range("E25").Value = '  +1.76%  '
$ws.Range("D26").Value = '158.12'
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("D27").Value = '16.57'
$ws.Range("E27").Value = '  +3.93%  '
$ws.Range("D28").Value = '7.18'
$ws.Range("E28").Value = '  +7.31%  '
$ws.Range("E29").Value = '  +2.18%  '
$ws.Range("D30").Value = '0.998'
$ws.Range("E30").Value = '  -0.26%  '
$ws.Range("D31").Value = '3.88'
$ws.Range("E31").Value = '  +11.33%  '
$ws.Range("D32").Value = '0.0514'
$ws.Range("E32").Value = '  +3.55%  '
$ws.Range("E33").Value = '  +4.68%  '
$ws.Range("D34").Value = '3.53'
$ws.Range("E34").Value = '  +7.06%  '
$ws.Range("D35").Value = '1.545.52'
$ws.Range("E35").Value = '  +2.33%  '
$ws.Range("E36").Value = '  +2.64%  '
$ws.Range("D37").Value = '1.09'
$ws.Range("E37").Value = '  +6.14%  '
$ws.Range("D38").Value = '85.18'
$ws.Range("E38").Value = '  +1.43%  '
$ws.Range("D39").Value = '0.0188'
$ws.Range("E39").Value = '  +5.09%  '
$ws.Range("D40").Value = '0.625'
$ws.Range("E40").Value = '  +3.73%  '
$ws.Range("D41").Value = '2.80'
$ws.Range("E41").Value = '  +3.68%  '
$ws.Range("E42").Value = '  +1.46%  '
$ws.Range("E43").Value = '  +8.95%  '
$ws.Range("D44").Value = '2.16'
$ws.Range("E44").Value = '  +8.56%  '
$ws.Range("E45").Value = '  +4.56%  '
$ws.Range("E46").Value = '  +4.40%  '
$ws.Range("D47").Value = '1.980.18'
$ws.Range("E47").Value = '  +9.17%  '
$ws.Range("D48").Value = '5.74'
$ws.Range("E48").Value = '  +2.27%  '
$ws.Range("E49").Value = '  -0.20%  '
$ws.Range("D50").Value = '52.70'
$ws.Range("E50").Value = '  +2.50%  '
$ws.Range("D51").Value = '11.83'
$ws.Range("E51").Value = '  +21.18%  '

# Restore default styling on the Price column (remove the temporary text
# number format) now that the values are committed as text.
$ws.Range("D2:D51").Style = "Normal"
